# Add new wheel-related parameters as two separate rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "W_WheelCenter"
$ws.Range("B8").Value = 20
$ws.Range("C8").Value = "mm"

$ws.Range("A9").Value = "T_WheelToShaft"
$ws.Range("B9").Value = 0.2
$ws.Range("C9").Value = "mm"
$ws.Range("D9").Value = "Distance between shaft and wheel"

$ws.Range("B9").Select()
